$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (Price) from Excel auto-converting numeric-looking
# text into floating point numbers, so exact string formatting
# (trailing zeros, thousand-dot separators, etc.) is preserved.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '29.042.89'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '1.830.19'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('D4').Value = '0.9988'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '241.32'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').Value = '0.6235'
$ws.Range('E6').Value = '  -5.43%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.07560'
$ws.Range('E8').Value = '  +2.06%  '
$ws.Range('E9').Value = '  +6.61%  '
$ws.Range('D10').Value = '0.2912'
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('D11').Value = '22.81'
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '0.07635'
$ws.Range('E12').Value = '  -1.77%  '
$ws.Range('D13').Value = '1.826.30'
$ws.Range('E13').Value = '  +0.65%  '
$ws.Range('D14').Value = '4.961'
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('D15').Value = '0.6651'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').Value = '82.35'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').Value = '0.000009060'
$ws.Range('E17').Value = '  +7.79%  '
$ws.Range('D18').Value = '6.006'
$ws.Range('E18').Value = '  -1.65%  '
$ws.Range('D19').Value = '29.039.77'
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '224.86'
$ws.Range('E20').Value = '  -1.14%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = '12.34'
$ws.Range('E21').Value = '  -0.94%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '7.195'
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '159.75'
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '8.393'
$ws.Range('E26').Value = '  -2.49%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = '0.1358'
$ws.Range('E27').Value = '  -2.39%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '17.84'
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = '1.494'
$ws.Range('E29').Value = '  -1.75%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '4.053'
$ws.Range('E30').Value = '  -1.48%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = '1.204'
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '4.025'
$ws.Range('E32').Value = '  -0.59%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.05209'
$ws.Range('E33').Value = '  -1.16%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '1.839'
$ws.Range('E34').Value = '  -1.34%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.155'
$ws.Range('E35').Value = '  +1.18%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.7323'
$ws.Range('E36').Value = '  -1.37%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = '2.610'
$ws.Range('E37').Value = '  -1.56%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.278.88'
$ws.Range('E38').Value = '  -1.80%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.749'
$ws.Range('E39').Value = '  +0.57%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.01779'
$ws.Range('E40').Value = '  -0.69%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '6.363'
$ws.Range('E41').Value = '  +7.40%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.8907'
$ws.Range('E42').Value = '  -3.99%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '101.43'
$ws.Range('E44').Value = '  -0.98%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.979.76'
$ws.Range('E45').Value = '  +1.05%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').Value = '0.5116'
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '63.46'
$ws.Range('E47').Value = '  +0.86%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.00000000119'
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').Value = '0.3965'
$ws.Range('E49').Value = '  -1.15%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '8.868'
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '1.655'
$ws.Range('E51').Value = '  -5.42%  '

# Remove the temporary text-format style again so the cells end up
# with no explicit style, matching the original (unstyled) cells.
$priceRange.ClearFormats()
